# "Fruta / hortaliza, semanal"
#
# The sheet holds one data row per "week" of Coliflor prices (rows 2..273,
# header on row 1). This edit inserts a brand-new most-recent observation
# at the top of the data block (row 158) and pushes every subsequent
# observation (rows 158..273) down by one row, spilling the previously
# last row (273) into a newly created row 274.
#
# Only the "observation" columns shift: D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# O (Origen) and P (Precio $/Kg). The descriptive columns A, B, C, E, F,
# G, H, N, Q, R are constant for every row in this block, so they do not
# need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

# --- Step 1: materialise the new last row (274) as an exact copy of the
#     current last row (273), BEFORE anything is shifted. ---
$fmt = $ws.Range("D273").NumberFormat()
$ws.Range("D274").NumberFormat = $fmt

$srcRow = $ws.Range("A273:R273")
$dstRow = $ws.Range("A274:R274")
$dstRow.Value = $srcRow.Value()

# --- Step 2: ripple every row down by one: row r (159..273) takes on the
#     "observation" values that row (r-1) held before this script ran.
#     Walking from the bottom up means each source row is read before it
#     is itself overwritten. ---
for ($r = 273; $r -ge 159; $r--) {
    $prev = $r - 1
    foreach ($col in $cols) {
        $cell = $col + $prev
        $target = $col + $r
        $ws.Range($target).Value = $ws.Range($cell).Value()
    }
}

# --- Step 3: row 158 becomes the brand-new observation (not sourced from
#     any existing row). Calidad (I) and Origen (O) stay "Primera" /
#     "Region del Maule" as they already were. ---
$ws.Range("D158").Value = 44762
$ws.Range("J158").Value = 2500
$ws.Range("K158").Value = 1500
$ws.Range("L158").Value = 1500
$ws.Range("M158").Value = 1500
$ws.Range("P158").Value = 1500
